$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2-39 to reflect repulled data / recalculated means
$dsfValues = @{
    2 = -1
    3 = -1
    4 = 7
    5 = 4
    6 = -2
    7 = 1
    8 = -2
    9 = -4
    10 = -4
    11 = 1
    12 = 1
    13 = 5
    14 = 5
    15 = 5
    16 = -3
    17 = 0
    18 = 7
    19 = -3
    20 = 6
    21 = 3
    22 = 7
    23 = 4
    24 = 1
    25 = 4
    26 = 0
    27 = -2
    28 = -2
    29 = 1
    30 = 7
    31 = 0
    32 = 3
    33 = -2
    34 = 6
    35 = 0
    36 = -2
    37 = 0
    38 = -3
    39 = 0
}

foreach ($row in $dsfValues.Keys | Sort-Object) {
    $ws.Cells.Item($row, 6).Value = $dsfValues[$row]
}

